# fastexcel "fill/byName.xlsx" test-resource refresh (#249):
# Sheet1's fill-template demo is translated from Chinese labels/placeholders
# to English ones. Sheet2/Sheet3 keep their original content untouched.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Header row (row 1): translate the column labels.
$ws1.Range("A1").Value = "Name"
$ws1.Range("B1").Value = "Number"

# Data row (row 2): translate the template text for the "complex"/"ignored"
# columns; A2/B2 ({name}/{number}) stay as-is.
$ws1.Range("C2").Value = "{name} is {number} years old this year"
$ws1.Range("D2").Value = "\{name\} ignored, {name}"

$ws1.Range("E1").Value = "Empty"
$ws1.Range("E2").Value = "Empty{.empty}"

$ws1.Range("D1").Value = "Ignored"
$ws1.Range("C1").Value = "Complex"

# These cells now hold Latin text, so re-apply the (CJK-capable) font to them
# the way Excel does when you retype a cell's content.
$ws1.Range("A1:D1").Font.Name = "宋体"
$ws1.Range("C2:D2").Font.Name = "宋体"
$ws1.Range("E1:E2").Font.Name = "宋体"

# Leave the cursor parked on C2, matching the saved selection.
[void]$ws1.Range("C2").Select()
